$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column K (reuse the same formatting as the other header cells,
# e.g. A1, by copying its format onto K1)
$ws.Range("K1").Value = "intervention_type"
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# Values for K2:K15
$values = @{
    2  = "BEHAVIORAL"
    3  = "PROCEDURE"
    4  = "BEHAVIORAL"
    5  = "OTHER"
    6  = "OTHER"
    7  = "OTHER"
    8  = "DRUG"
    9  = "OTHER"
    10 = "DRUG"
    11 = "DRUG"
    12 = "DEVICE"
    13 = "DRUG"
    14 = "RADIATION"
    15 = "OTHER"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 11).Value = $values[$row]
}

# K16 left empty (matches existing pattern of blank inline string cells like C16, E16)
# Touch the cell's formatting (without assigning a value) so it remains present
# in the sheet as an empty cell rather than being dropped entirely.
$ws.Range("K16").Borders.LineStyle = -4142
